# Screenshots for gameweek 18 — append 9 new match rows (157-165) to the
# results table and extend the "100 - possession" helper formula down to
# the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Home, Away, Home xG, Away xG, Home Possession
$rows = @(
    @{ Row = 157; Home = "Everton";                 Away = "West Ham";                 HomeXG = "0.371045"; AwayXG = "1.19713"; Poss = 60 },
    @{ Row = 158; Home = "Manchester United";        Away = "Aston Villa";              HomeXG = "2.46473";  AwayXG = "1.56434"; Poss = 48.4 },
    @{ Row = 159; Home = "Tottenham";                Away = "Leeds";                    HomeXG = "2.5286";   AwayXG = "1.15626"; Poss = 36.7 },
    @{ Row = 160; Home = "Crystal Palace";            Away = "Sheffield United";        HomeXG = "0.598888"; AwayXG = "0.323158"; Poss = 44.4 },
    @{ Row = 161; Home = "Brighton";                  Away = "Wolverhampton Wanderers"; HomeXG = "2.4228";   AwayXG = "1.81856"; Poss = 54.5 },
    @{ Row = 162; Home = "West Bromwich Albion";      Away = "Arsenal";                 HomeXG = "0.869874"; AwayXG = "4.03429"; Poss = 38.1 },
    @{ Row = 163; Home = "Newcastle United";          Away = "Leicester";               HomeXG = "0.264584"; AwayXG = "0.702081"; Poss = 43.3 },
    @{ Row = 164; Home = "Chelsea";                   Away = "Manchester City";         HomeXG = "0.691013"; AwayXG = "3.18063"; Poss = 54.6 },
    @{ Row = 165; Home = "Southampton";               Away = "Liverpool";               HomeXG = "0.500436"; AwayXG = "1.37204"; Poss = 33 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = $r.Home
    $ws.Range("B$rowNum").Value = $r.Away
    $ws.Range("E$rowNum").Value = $r.Poss
}

# Columns C/D hold the xG figures as TEXT (shared strings) in the source
# workbook, not numbers — force text storage with a quote-prefix, then
# strip the resulting number-format style back off so no stray formatting
# is left behind on the cells.
$cdRange = $ws.Range("C157:D165")
$origStyle = $ws.Range("A1").Style

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("C$rowNum").Value = "'" + $r.HomeXG
    $ws.Range("D$rowNum").Value = "'" + $r.AwayXG
}
$cdRange.Style = $origStyle

# Extend the shared "100 - possession" formula down through the new rows.
$ws.Range("F157:F165").Formula = "=SUM(100-E157)"

# Refresh the used-range dimension / selection to mirror the saved file.
$ws.Range("E166").Select()
